$d = $word.ActiveDocument

# --- Paragraph 3: drop the proofErr (gramStart/gramEnd) markers around "let's" ---
# Runs stay split exactly as they were; only the two <w:proofErr/> tags disappear.
# Target the range up to (but excluding) the paragraph's own end-of-paragraph mark so
# InsertXML replaces the paragraph's run content in place instead of splitting it.
$p1 = $d.Paragraphs(3)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="4892DF69" w14:textId="36741972" w:rsidR="00F46537" w:rsidRDefault="00F46537" w:rsidP="00F46537">' +
        '<w:r><w:t xml:space="preserve">This is extremely hard to complete, 90% of people will fail this task, </w:t></w:r>' +
        '<w:r w:rsidR="00631380"><w:t>let’s</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> see how smart you are</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Paragraph 11: drop the proofErr (spellStart/spellEnd) markers and fix "tis" -> "this",
#     split across three runs: "...making it t" / "h" / "is far" ---
# This is the last paragraph in the body (right before sectPr), so again stop one
# character short of its Range.End to avoid touching the end-of-paragraph mark.
$p2 = $d.Paragraphs(11)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="0E38216C" w14:textId="232C87DB" w:rsidR="00F46537" w:rsidRDefault="00F46537" w:rsidP="00F46537">' +
        '<w:r><w:t>CONGRATS, you completed Mastermind, you are a vary smart person for making it t</w:t></w:r>' +
        '<w:r><w:t>h</w:t></w:r>' +
        '<w:r><w:t>is far</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)
